# Auto-generated edit script for table_es1b.xlsx
# Commit: 2017-01-31 update: energy.gov - chunk 7
# Updates the workbook from "October 2016 YTD" reporting period to "November 2016 YTD"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Title / header text updates (January through October -> November) ---
$ws.Range("A3").Value = "Net Generation and Consumption of Fuels for January through November"
$ws.Range("C6").Value = "November 2016 YTD"
$ws.Range("D6").Value = "November 2015 YTD"
$ws.Range("F6").Value = "November 2016 YTD"
$ws.Range("G6").Value = "November 2015 YTD"
$ws.Range("H6").Value = "November 2016 YTD"
$ws.Range("I6").Value = "November 2015 YTD"
$ws.Range("J6").Value = "November 2016 YTD"
$ws.Range("K6").Value = "November 2015 YTD"
$ws.Range("L6").Value = "November 2016 YTD"
$ws.Range("M6").Value = "November 2015 YTD"
$ws.Range("N6").Value = "November 2016 YTD"
$ws.Range("O6").Value = "November 2015 YTD"
$ws.Range("A45").Value = "Sales, Revenue, and Average Price of Electricity to Ultimate Customers for January through November"
$ws.Range("B48").Value = "November 2016 YTD"
$ws.Range("C48").Value = "November 2015 YTD"
$ws.Range("E48").Value = "November 2016 YTD"
$ws.Range("F48").Value = "November 2015 YTD"
$ws.Range("H48").Value = "November 2016 YTD"
$ws.Range("I48").Value = "November 2015 YTD"

# --- Data table updates: Net Generation (Thousand Megawatthours), rows 8-23 ---
# Row 8
$ws.Range("C8").Value = 1121120
$ws.Range("D8").Value = 1262903
$ws.Range("E8").Value = -0.112
$ws.Range("F8").Value = 835856
$ws.Range("G8").Value = 929827
$ws.Range("H8").Value = 276328
$ws.Range("I8").Value = 322545
$ws.Range("J8").Value = 391
$ws.Range("K8").Value = 468
$ws.Range("L8").Value = 8544
$ws.Range("M8").Value = 10064

# Row 9
$ws.Range("C9").Value = 11528
$ws.Range("D9").Value = 16424
$ws.Range("E9").Value = -0.298
$ws.Range("F9").Value = 7913
$ws.Range("G9").Value = 9736
$ws.Range("H9").Value = 3060
$ws.Range("I9").Value = 5987
$ws.Range("J9").Value = 96
$ws.Range("K9").Value = 176
$ws.Range("L9").Value = 459
$ws.Range("M9").Value = 525

# Row 10
$ws.Range("C10").Value = 10366
$ws.Range("D10").Value = 10128
$ws.Range("E10").Value = 0.023
$ws.Range("F10").Value = 8214
$ws.Range("G10").Value = 7674
$ws.Range("H10").Value = 1276
$ws.Range("I10").Value = 1533
$ws.Range("J10").Value = 4
$ws.Range("K10").Value = 8
$ws.Range("L10").Value = 871
$ws.Range("M10").Value = 912

# Row 11
$ws.Range("C11").Value = 1284457
$ws.Range("D11").Value = 1223705
$ws.Range("E11").Value = 0.05
$ws.Range("F11").Value = 608453
$ws.Range("G11").Value = 565471
$ws.Range("H11").Value = 584558
$ws.Range("I11").Value = 571163
$ws.Range("J11").Value = 7148
$ws.Range("K11").Value = 6853
$ws.Range("L11").Value = 84298
$ws.Range("M11").Value = 80218

# Row 12
$ws.Range("C12").Value = 11987
$ws.Range("D12").Value = 12007
$ws.Range("E12").Value = -0.002
$ws.Range("F12").Value = 142
$ws.Range("G12").Value = 197
$ws.Range("H12").Value = 3598
$ws.Range("I12").Value = 3214
$ws.Range("L12").Value = 8248
$ws.Range("M12").Value = 8595

# Row 13
$ws.Range("C13").Value = 733632
$ws.Range("D13").Value = 727544
$ws.Range("E13").Value = 0.008
$ws.Range("F13").Value = 387127
$ws.Range("G13").Value = 380683
$ws.Range("H13").Value = 346505
$ws.Range("I13").Value = 346861

# Row 14
$ws.Range("C14").Value = 243220
$ws.Range("D14").Value = 225915
$ws.Range("E14").Value = 0.077
$ws.Range("F14").Value = 225608
$ws.Range("G14").Value = 208344
$ws.Range("H14").Value = 16381
$ws.Range("I14").Value = 16275
$ws.Range("J14").Value = 54
$ws.Range("K14").Value = 31
$ws.Range("L14").Value = 1177
$ws.Range("M14").Value = 1265

# Row 15
$ws.Range("C15").Value = 310064
$ws.Range("D15").Value = 266559
$ws.Range("E15").Value = 0.163
$ws.Range("F15").Value = 38277
$ws.Range("G15").Value = 33811
$ws.Range("H15").Value = 243087
$ws.Range("I15").Value = 203619
$ws.Range("J15").Value = 2921
$ws.Range("K15").Value = 2954
$ws.Range("L15").Value = 25779
$ws.Range("M15").Value = 26175

# Row 16
$ws.Range("C16").Value = 203453
$ws.Range("D16").Value = 170620
$ws.Range("E16").Value = 0.192
$ws.Range("F16").Value = 31209
$ws.Range("G16").Value = 27308
$ws.Range("H16").Value = 172047
$ws.Range("I16").Value = 143159
$ws.Range("J16").Value = 128
$ws.Range("K16").Value = 106
$ws.Range("L16").Value = 69
$ws.Range("M16").Value = 47

# Row 17
$ws.Range("C17").Value = 33832
$ws.Range("D17").Value = 23323
$ws.Range("E17").Value = 0.451
$ws.Range("F17").Value = 2033
$ws.Range("G17").Value = 1396
$ws.Range("H17").Value = 31237
$ws.Range("I17").Value = 21514
$ws.Range("J17").Value = 534
$ws.Range("K17").Value = 393
$ws.Range("L17").Value = 28
$ws.Range("M17").Value = 20

# Row 18
$ws.Range("C18").Value = 36842
$ws.Range("D18").Value = 38341
$ws.Range("E18").Value = -0.039
$ws.Range("F18").Value = 2697
$ws.Range("G18").Value = 2759
$ws.Range("H18").Value = 9487
$ws.Range("I18").Value = 10550
$ws.Range("J18").Value = 65
$ws.Range("K18").Value = 45
$ws.Range("L18").Value = 24593
$ws.Range("M18").Value = 24987

# Row 19
$ws.Range("C19").Value = 20140
$ws.Range("D19").Value = 19734
$ws.Range("E19").Value = 0.021
$ws.Range("F19").Value = 1341
$ws.Range("G19").Value = 1357
$ws.Range("H19").Value = 15517
$ws.Range("I19").Value = 14846
$ws.Range("J19").Value = 2193
$ws.Range("K19").Value = 2411
$ws.Range("L19").Value = 1089
$ws.Range("M19").Value = 1121

# Row 20
$ws.Range("C20").Value = 15797
$ws.Range("D20").Value = 14541
$ws.Range("E20").Value = 0.086
$ws.Range("F20").Value = 997
$ws.Range("G20").Value = 991
$ws.Range("H20").Value = 14800
$ws.Range("I20").Value = 13549

# Row 21
$ws.Range("C21").Value = -5933
$ws.Range("D21").Value = -4811
$ws.Range("E21").Value = 0.233
$ws.Range("F21").Value = -4972
$ws.Range("G21").Value = -3895
$ws.Range("H21").Value = -961
$ws.Range("I21").Value = -916

# Row 22
$ws.Range("C22").Value = 12550
$ws.Range("D22").Value = 12800
$ws.Range("E22").Value = -0.019
$ws.Range("F22").Value = 288
$ws.Range("G22").Value = 509
$ws.Range("H22").Value = 6494
$ws.Range("I22").Value = 6231
$ws.Range("J22").Value = 993
$ws.Range("K22").Value = 1072
$ws.Range("L22").Value = 4776
$ws.Range("M22").Value = 4988

# Row 23
$ws.Range("C23").Value = 3732992
$ws.Range("D23").Value = 3753174
$ws.Range("F23").Value = 2106905
$ws.Range("G23").Value = 2132358
$ws.Range("H23").Value = 1480326
$ws.Range("I23").Value = 1476513
$ws.Range("J23").Value = 11608
$ws.Range("K23").Value = 11562
$ws.Range("L23").Value = 134153
$ws.Range("M23").Value = 132742

# Estimated Solar rows 25-27
# Row 25
$ws.Range("C25").Value = 18281
$ws.Range("D25").Value = 13225
$ws.Range("E25").Value = 0.382
$ws.Range("J25").Value = 6710
$ws.Range("K25").Value = 5340
$ws.Range("L25").Value = 1720
$ws.Range("M25").Value = 1358
$ws.Range("N25").Value = 9851
$ws.Range("O25").Value = 6527

# Row 26
$ws.Range("C26").Value = 48820
$ws.Range("D26").Value = 33447
$ws.Range("E26").Value = 0.46
$ws.Range("F26").Value = 1960
$ws.Range("G26").Value = 1292
$ws.Range("H26").Value = 28017
$ws.Range("I26").Value = 18517
$ws.Range("J26").Value = 7244
$ws.Range("K26").Value = 5733
$ws.Range("L26").Value = 1748
$ws.Range("M26").Value = 1378
$ws.Range("N26").Value = 9851
$ws.Range("O26").Value = 6527

# Row 27
$ws.Range("C27").Value = 52113
$ws.Range("D27").Value = 36548
$ws.Range("E27").Value = 0.426
$ws.Range("F27").Value = 2033
$ws.Range("G27").Value = 1396
$ws.Range("H27").Value = 31237
$ws.Range("I27").Value = 21514
$ws.Range("J27").Value = 7244
$ws.Range("K27").Value = 5733
$ws.Range("L27").Value = 1748
$ws.Range("M27").Value = 1378
$ws.Range("N27").Value = 9851
$ws.Range("O27").Value = 6527

# Consumption of Fossil Fuels for Electricity Generation, rows 29-32
# Row 29
$ws.Range("C29").Value = 613093
$ws.Range("D29").Value = 689370
$ws.Range("E29").Value = -0.111
$ws.Range("F29").Value = 450965
$ws.Range("G29").Value = 501628
$ws.Range("H29").Value = 158865
$ws.Range("I29").Value = 183886
$ws.Range("J29").Value = 133
$ws.Range("K29").Value = 149
$ws.Range("L29").Value = 3130
$ws.Range("M29").Value = 3707

# Row 30
$ws.Range("C30").Value = 19315
$ws.Range("D30").Value = 27344
$ws.Range("E30").Value = -0.294
$ws.Range("F30").Value = 14316
$ws.Range("G30").Value = 17385
$ws.Range("H30").Value = 4363
$ws.Range("I30").Value = 9119
$ws.Range("J30").Value = 119
$ws.Range("K30").Value = 241
$ws.Range("L30").Value = 517
$ws.Range("M30").Value = 599

# Row 31
$ws.Range("C31").Value = 3939
$ws.Range("D31").Value = 3768
$ws.Range("E31").Value = 0.045
$ws.Range("F31").Value = 3170
$ws.Range("G31").Value = 2888
$ws.Range("H31").Value = 544
$ws.Range("I31").Value = 643
$ws.Range("L31").Value = 224
$ws.Range("M31").Value = 235

# Row 32
$ws.Range("C32").Value = 9698844
$ws.Range("D32").Value = 9209356
$ws.Range("E32").Value = 0.053
$ws.Range("F32").Value = 4702672
$ws.Range("G32").Value = 4351897
$ws.Range("H32").Value = 4344715
$ws.Range("I32").Value = 4225560
$ws.Range("J32").Value = 63156
$ws.Range("K32").Value = 64344
$ws.Range("L32").Value = 588301
$ws.Range("M32").Value = 567555

# Consumption of Fossil Fuels for Useful Thermal Output, rows 34-37
# Row 34
$ws.Range("C34").Value = 13016
$ws.Range("D34").Value = 15269
$ws.Range("E34").Value = -0.148
$ws.Range("F34").Value = 919
$ws.Range("G34").Value = 937
$ws.Range("H34").Value = 1594
$ws.Range("I34").Value = 1829
$ws.Range("J34").Value = 488
$ws.Range("K34").Value = 577
$ws.Range("L34").Value = 10016
$ws.Range("M34").Value = 11927

# Row 35
$ws.Range("C35").Value = 2307
$ws.Range("D35").Value = 2932
$ws.Range("E35").Value = -0.213
$ws.Range("G35").Value = 61
$ws.Range("H35").Value = 919
$ws.Range("I35").Value = 1065
$ws.Range("J35").Value = 113
$ws.Range("K35").Value = 277
$ws.Range("L35").Value = 1259
$ws.Range("M35").Value = 1529

# Row 36
$ws.Range("C36").Value = 923
$ws.Range("D36").Value = 1058
$ws.Range("E36").Value = -0.128
$ws.Range("G36").Value = 9
$ws.Range("H36").Value = 95
$ws.Range("I36").Value = 99
$ws.Range("K36").Value = 15
$ws.Range("L36").Value = 818
$ws.Range("M36").Value = 935

# Row 37
$ws.Range("C37").Value = 876439
$ws.Range("D37").Value = 853730
$ws.Range("E37").Value = 0.027
$ws.Range("F37").Value = 10063
$ws.Range("G37").Value = 7317
$ws.Range("H37").Value = 262056
$ws.Range("I37").Value = 259989
$ws.Range("J37").Value = 43852
$ws.Range("K37").Value = 42266
$ws.Range("L37").Value = 560468
$ws.Range("M37").Value = 544157

# Consumption of Fossil Fuels for Electricity Generation and Useful Thermal Output, rows 39-42
# Row 39
$ws.Range("C39").Value = 626108
$ws.Range("D39").Value = 704639
$ws.Range("E39").Value = -0.111
$ws.Range("F39").Value = 451884
$ws.Range("G39").Value = 502565
$ws.Range("H39").Value = 160458
$ws.Range("I39").Value = 185714
$ws.Range("J39").Value = 621
$ws.Range("K39").Value = 726
$ws.Range("L39").Value = 13145
$ws.Range("M39").Value = 15634

# Row 40
$ws.Range("C40").Value = 21622
$ws.Range("D40").Value = 30277
$ws.Range("E40").Value = -0.286
$ws.Range("F40").Value = 14333
$ws.Range("G40").Value = 17446
$ws.Range("H40").Value = 5281
$ws.Range("I40").Value = 10184
$ws.Range("J40").Value = 232
$ws.Range("K40").Value = 518
$ws.Range("L40").Value = 1776
$ws.Range("M40").Value = 2128

# Row 41
$ws.Range("C41").Value = 4862
$ws.Range("D41").Value = 4826
$ws.Range("E41").Value = 0.007
$ws.Range("F41").Value = 3172
$ws.Range("G41").Value = 2896
$ws.Range("H41").Value = 639
$ws.Range("I41").Value = 742
$ws.Range("J41").Value = 9
$ws.Range("K41").Value = 17
$ws.Range("L41").Value = 1042
$ws.Range("M41").Value = 1171

# Row 42
$ws.Range("C42").Value = 10575283
$ws.Range("D42").Value = 10063086
$ws.Range("E42").Value = 0.051
$ws.Range("F42").Value = 4712735
$ws.Range("G42").Value = 4359214
$ws.Range("H42").Value = 4606771
$ws.Range("I42").Value = 4485550
$ws.Range("J42").Value = 107008
$ws.Range("K42").Value = 106611
$ws.Range("L42").Value = 1148769
$ws.Range("M42").Value = 1111712

# Sales, Revenue, and Average Price of Electricity to Ultimate Customers, rows 49-53
# Row 49
$ws.Range("B49").Value = 1286581
$ws.Range("C49").Value = 1292426
$ws.Range("E49").Value = 161836
$ws.Range("F49").Value = 163864
$ws.Range("G49").Value = -0.012
$ws.Range("H49").Value = 12.58
$ws.Range("J49").Value = -0.008

# Row 50
$ws.Range("B50").Value = 1250159
$ws.Range("C50").Value = 1253922
$ws.Range("D50").Value = -0.003
$ws.Range("E50").Value = 129903
$ws.Range("F50").Value = 133956
$ws.Range("G50").Value = -0.03
$ws.Range("H50").Value = 10.39
$ws.Range("I50").Value = 10.68
$ws.Range("J50").Value = -0.027

# Row 51
$ws.Range("B51").Value = 861283
$ws.Range("C51").Value = 908283
$ws.Range("D51").Value = -0.052
$ws.Range("E51").Value = 58213
$ws.Range("F51").Value = 63123
$ws.Range("G51").Value = -0.078
$ws.Range("H51").Value = 6.76
$ws.Range("I51").Value = 6.95
$ws.Range("J51").Value = -0.027

# Row 52
$ws.Range("B52").Value = 6845
$ws.Range("C52").Value = 7017
$ws.Range("E52").Value = 649
$ws.Range("F52").Value = 710
$ws.Range("H52").Value = 9.49
$ws.Range("I52").Value = 10.11
$ws.Range("J52").Value = -0.061

# Row 53
$ws.Range("B53").Value = 3404868
$ws.Range("C53").Value = 3461649
$ws.Range("D53").Value = -0.016
$ws.Range("E53").Value = 350601
$ws.Range("F53").Value = 361654
$ws.Range("G53").Value = -0.031
$ws.Range("H53").Value = 10.3
$ws.Range("I53").Value = 10.45
$ws.Range("J53").Value = -0.014

